$wb = $excel.ActiveWorkbook

# --- "US 3" sheet: task status change from "ToDo" to "Completed" ---
$ws3 = $wb.Worksheets.Item("US 3")
$ws3.Range("F8").Value = "Completed"

# --- "US 5" sheet: remaining-time re-estimate + removal of completed/duplicate task row ---
$ws5 = $wb.Worksheets.Item("US 5")
$ws5.Range("D5").Value = 3
$ws5.Rows("6").Delete()
$ws5.Range("D5").Select()

# --- restore per-sheet active cell selections ---
$ws3.Range("F8").Select()

# --- "US7" sheet becomes the active/visible tab ---
$ws7 = $wb.Worksheets.Item("US7")
$ws7.Activate()
$ws7.Range("E22").Select()
